$wb = $excel.ActiveWorkbook

# --- "About" sheet: remove the stray date value in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Clear()

# --- "Set Values Here" sheet: update the "fuel tax revenue" row (row 9) weights ---
$wsValues = $wb.Worksheets.Item("Set Values Here")
$wsValues.Range("C9").Value = 5
$wsValues.Range("D9").Value = 0
$wsValues.Range("F9").Value = 5

# Update the selected cell on this sheet to match the authored state,
# then restore "About" as the active sheet/tab.
$wsValues.Range("F10").Select()
$wsAbout.Activate()
